$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# The underlying source data was re-aggregated: a handful of "Ar Condicionado"
# / "Materiais elétricos" rows were inserted into the per-building category
# breakdown (pushing every following row down), and the sum/count figures for
# several categories shifted accordingly. Easiest + most faithful way to
# reproduce that here is to rewrite the whole A2:D55 data block in place.
# ---------------------------------------------------------------------------

# Clear the previous data body first (old sheet only went down to row 50).
# Use ClearContents (not Clear) so the existing bold/border/centered style
# on columns A/B stays in place for the values we're about to write back in.
$ws.Range("A2:D50").ClearContents()

$csv = @"
2|A|Ar Condicionado|19.853712|132
3|A|Computadores|4.76096|44
4|A|Eletrodomésticos|1.882784|6
5|A|Eletrônicos|2.6255|48
6|A|Ferramentas|0.136128|2
7|A|Livros|0.12|1
8|A|Maquinário|0.06875000000000001|1
9|A|Materiais elétricos|0.109956|1
10|A|Material de Escritório|4.701928|42
11|A|Mobília|571.668886|1308
12|A|Outros|26.846936|102
13|A|Periféricos de Informática|0.544|5
14|A|Transporte Particular|0.034|3
15|A|Vidraria|0.301952|13
16|B|Ar Condicionado|16.152456|118
17|B|Computadores|13.607344|119
18|B|Eletrodomésticos|8.239568|25
19|B|Eletrônicos|26.32633|201
20|B|Ferramentas|2.273606|19
21|B|Livros|0.565824|13
22|B|Maquinário|93.317637|151
23|B|Maquinário Pesado|1.46412|3
24|B|Maquinário Sensível|0.274304|4
25|B|Materiais elétricos|0.817852|10
26|B|Material de Escritório|7.047371|53
27|B|Mobília|555.623498|1200
28|B|Outros|731.8399280000001|475
29|B|Periféricos de Informática|1.912935|21
30|B|Reagentes|10.506306|179
31|B|Transporte Particular|2.514475|19
32|B|Vidraria|15.324944|148
33|C|Ar Condicionado|9.177708000000001|64
34|C|Computadores|12.453464|104
35|C|Eletrodomésticos|0.24|2
36|C|Eletrônicos|12.698566|210
37|C|Ferramentas|0.06467400000000001|2
38|C|Livros|0.24|2
39|C|Maquinário|0.25428|5
40|C|Materiais elétricos|0.17707|4
41|C|Material de Escritório|109.565079|274
42|C|Mobília|205.244884|384
43|C|Outros|11.979093|58
44|C|Periféricos de Informática|1.898213|20
45|C|Transporte Particular|0.12|1
46|C|Vidraria|0.62|3
47|C|Ar Condicionado|0.12492|2
48|E|Eletrônicos|0.03|1
49|E|Ferramentas|22.928|76
50|E|Materiais elétricos|2.64|22
51|E|Mobília|23.870035|37
52|S|Outros|0.128|1
53|V|Ar Condicionado|0.54208|4
54|V|Computadores|6.96|58
55|V|Mobília|36.195091|115
"@

$lines = $csv -split "`r?`n"
$lastRow = 1
foreach ($line in $lines) {
    if ($line.Trim().Length -eq 0) { continue }
    $parts = $line -split '\|'
    $r = [int]$parts[0]
    $ws.Cells.Item($r, 1).Value = $parts[1]
    $ws.Cells.Item($r, 2).Value = $parts[2]
    $ws.Cells.Item($r, 3).Value = [double]$parts[3]
    $ws.Cells.Item($r, 4).Value = [int]$parts[4]
    if ($r -gt $lastRow) { $lastRow = $r }
}

# Rows beyond the workbook's original extent (row 50) don't inherit the
# bold/bordered/centered style used by columns A and B in the data rows
# above, so copy that formatting down explicitly.
if ($lastRow -gt 50) {
    $ws.Range("A2:B2").Copy() | Out-Null
    $ws.Range("A51:B" + $lastRow).PasteSpecial(-4122) | Out-Null
    $excel.CutCopyMode = 0
}

Write-Output "done"
